# Normalize accented Portuguese text in NOME_UNIDADE (col B) and
# NIVEL_CURSO (col C) columns: strip diacritics (cedilla / tilde / acute)
# so that e.g. "Dança" -> "Danca", "Educação Física" -> "Educacao Fisica",
# "Graduação" -> "Graduacao".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value()
    switch ($bVal) {
        "Dança - Licenciatura"          { $bCell.Value = "Danca - Licenciatura" }
        "Educação Física"               { $bCell.Value = "Educacao Fisica" }
        "Educação Física - Bacharelado" { $bCell.Value = "Educacao Fisica - Bacharelado" }
    }

    $cCell = $ws.Cells.Item($r, 3)
    $cVal = $cCell.Value()
    if ($cVal -eq "Graduação") {
        $cCell.Value = "Graduacao"
    }
}
